{"js": "const newValues = [\"63-25=38\", \"25+39=64\", \"31-17=14\", \"50-14=36\", \"48+16=64\", \"96-57=39\", \"9+8=17\", \"62-19=43\", \"25+29=54\", \"28+3=31\", \"39+46=85\", \"19+32=51\", \"70-1=69\", \"78+4=82\", \"71-3=68\", \"81-36=45\", \"49+12=61\", \"54-16=38\", \"50-35=15\", \"72-4=68\", \"94-18=76\", \"9+36=45\", \"69+19=88\", \"48+47=95\", \"95-78=17\", \"46+9=55\", \"91-49=42\", \"61-27=34\", \"18+79=97\", \"91-38=53\", \"92-33=59\", \"32+49=81\", \"34+58=92\", \"17+14=31\", \"34-9=25\", \"27+55=82\", \"57+27=84\", \"71-55=16\", \"73-54=19\", \"86-68=18\", \"18+29=47\", \"57+39=96\", \"85-79=6\", \"74-56=18\", \"66+18=84\", \"6+5=11\", \"70-29=41\", \"29+19=48\", \"41-7=34\", \"38+28=66\", \"16+25=41\", \"56-9=47\", \"72-64=8\", \"62-8=54\", \"50-31=19\", \"71-12=59\", \"56+9=65\", \"91-23=68\", \"61-4=57\", \"34-8=26\", \"95-6=89\", \"37+57=94\", \"36+8=44\", \"7+64=71\", \"50-1=49\", \"38+56=94\", \"8+53=61\", \"18+19=37\", \"8+24=32\", \"65-26=39\", \"60-33=27\", \"92-73=19\", \"60-51=9\", \"49+12=61\", \"69+12=81\", \"59+6=65\", \"74-55=19\", \"37+49=86\", \"25+47=72\", \"85+7=92\", \"40-11=29\", \"64+17=81\", \"60-26=34\", \"7+85=92\", \"47+6=53\", \"74+19=93\", \"73-7=66\", \"28+26=54\", \"47+34=81\", \"28+35=63\", \"18+73=91\", \"28+54=82\", \"43+49=92\", \"91-26=65\", \"20-2=18\", \"56-9=47\", \"80-68=12\", \"74-68=6\", \"70-19=51\", \"7+47=54\"];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nlet idx = 0;\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n  for (const cell of cells.items) {\n    if (idx < newValues.length) {\n      cell.value = newValues[idx];\n      idx++;\n    }\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n$values = @(\"63-25=38\",\"25+39=64\",\"31-17=14\",\"50-14=36\",\"48+16=64\",\"96-57=39\",\"9+8=17\",\"62-19=43\",\"25+29=54\",\"28+3=31\",\"39+46=85\",\"19+32=51\",\"70-1=69\",\"78+4=82\",\"71-3=68\",\"81-36=45\",\"49+12=61\",\"54-16=38\",\"50-35=15\",\"72-4=68\",\"94-18=76\",\"9+36=45\",\"69+19=88\",\"48+47=95\",\"95-78=17\",\"46+9=55\",\"91-49=42\",\"61-27=34\",\"18+79=97\",\"91-38=53\",\"92-33=59\",\"32+49=81\",\"34+58=92\",\"17+14=31\",\"34-9=25\",\"27+55=82\",\"57+27=84\",\"71-55=16\",\"73-54=19\",\"86-68=18\",\"18+29=47\",\"57+39=96\",\"85-79=6\",\"74-56=18\",\"66+18=84\",\"6+5=11\",\"70-29=41\",\"29+19=48\",\"41-7=34\",\"38+28=66\",\"16+25=41\",\"56-9=47\",\"72-64=8\",\"62-8=54\",\"50-31=19\",\"71-12=59\",\"56+9=65\",\"91-23=68\",\"61-4=57\",\"34-8=26\",\"95-6=89\",\"37+57=94\",\"36+8=44\",\"7+64=71\",\"50-1=49\",\"38+56=94\",\"8+53=61\",\"18+19=37\",\"8+24=32\",\"65-26=39\",\"60-33=27\",\"92-73=19\",\"60-51=9\",\"49+12=61\",\"69+12=81\",\"59+6=65\",\"74-55=19\",\"37+49=86\",\"25+47=72\",\"85+7=92\",\"40-11=29\",\"64+17=81\",\"60-26=34\",\"7+85=92\",\"47+6=53\",\"74+19=93\",\"73-7=66\",\"28+26=54\",\"47+34=81\",\"28+35=63\",\"18+73=91\",\"28+54=82\",\"43+49=92\",\"91-26=65\",\"20-2=18\",\"56-9=47\",\"80-68=12\",\"74-68=6\",\"70-19=51\",\"7+47=54\")\n$rows = $table.Rows.Count\n$cols = $table.Columns.Count\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $table.Cell($r, $c)\n        $cell.Range.Text = $values[$idx]\n        $idx++\n    }\n}\n"}
